# Update Sheet1 cell values (A1:J20) per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 20,10

$data[0,0] = 148.229717625412
$data[0,1] = 111.604122310693
$data[0,2] = 13.2234045365934
$data[0,3] = 32.5750158320065
$data[0,4] = 118.067045751059
$data[0,5] = 157.447280063036
$data[0,6] = 182.540780763394
$data[0,7] = 80.4811589794611
$data[0,8] = 172.439766289871
$data[0,9] = 165.571800044538
$data[1,0] = 112.433960154855
$data[1,1] = 103.396801139878
$data[1,2] = 141.047907779481
$data[1,3] = 98.9095377265054
$data[1,4] = 163.981578016645
$data[1,5] = 100.985078094986
$data[1,6] = 56.6176130700007
$data[1,7] = 64.1690493860138
$data[1,8] = 128.105025704999
$data[1,9] = 31.6074366828461
$data[2,0] = 198.738175257453
$data[2,1] = 153.875191488245
$data[2,2] = 199.13598578383
$data[2,3] = 193.675257542019
$data[2,4] = 164.058820234639
$data[2,5] = 163.385138829884
$data[2,6] = 100.500519527355
$data[2,7] = 120.475367047114
$data[2,8] = 171.336902292137
$data[2,9] = 91.3067922421297
$data[3,0] = 6.58374904030177
$data[3,1] = 58.4566431392248
$data[3,2] = 95.5408974995561
$data[3,3] = 67.6215193549271
$data[3,4] = 13.1260387660591
$data[3,5] = 122.829748840458
$data[3,6] = 71.961750123632
$data[3,7] = 172.44797384946
$data[3,8] = 73.8501658075723
$data[3,9] = 133.056453118593
$data[4,0] = 129.587843143189
$data[4,1] = 95.1301574218693
$data[4,2] = 162.057416589026
$data[4,3] = 139.110995800752
$data[4,4] = 1.66322710069978
$data[4,5] = 5.37442770105527
$data[4,6] = 112.940371275386
$data[4,7] = 140.08451920938
$data[4,8] = 106.081263025329
$data[4,9] = 198.192902281039
$data[5,0] = 149.330897791931
$data[5,1] = 76.3430363854128
$data[5,2] = 145.286550347361
$data[5,3] = 140.395605256965
$data[5,4] = 179.028062233249
$data[5,5] = 194.354526137167
$data[5,6] = 112.468136526862
$data[5,7] = 19.5481469945741
$data[5,8] = 68.5161955973675
$data[5,9] = 154.681906921175
$data[6,0] = 56.946760535681
$data[6,1] = 62.0654137162796
$data[6,2] = 109.144256687325
$data[6,3] = 81.1329740477414
$data[6,4] = 158.988051004237
$data[6,5] = 53.9773170156299
$data[6,6] = 7.85590364032234
$data[6,7] = 73.4263884245541
$data[6,8] = 85.7834989604463
$data[6,9] = 41.1518291761874
$data[7,0] = 29.0233279713538
$data[7,1] = 84.1696392205402
$data[7,2] = 190.318883578442
$data[7,3] = 195.048572586406
$data[7,4] = 102.019593539657
$data[7,5] = 103.608017835584
$data[7,6] = 191.817774899219
$data[7,7] = 60.0249899830786
$data[7,8] = 192.01203044132
$data[7,9] = 158.684392533584
$data[8,0] = 50.4447675544977
$data[8,1] = 160.416000317976
$data[8,2] = 14.3941040217849
$data[8,3] = 173.144000011098
$data[8,4] = 141.975894450199
$data[8,5] = 130.240712654889
$data[8,6] = 113.170092791864
$data[8,7] = 155.145292242591
$data[8,8] = 88.5934571216784
$data[8,9] = 18.7715126288922
$data[9,0] = 10.3616123135954
$data[9,1] = 52.413603129058
$data[9,2] = 103.931778252093
$data[9,3] = 119.168258886397
$data[9,4] = 76.1096925829117
$data[9,5] = 67.5224294269096
$data[9,6] = 185.985900734545
$data[9,7] = 80.9244425412847
$data[9,8] = 180.122944796515
$data[9,9] = 147.68591008507
$data[10,0] = 197.518524060733
$data[10,1] = 39.5139828508319
$data[10,2] = 54.3010202489332
$data[10,3] = 64.9294338491417
$data[10,4] = 169.169574309685
$data[10,5] = 65.1612585713906
$data[10,6] = 86.0241528069713
$data[10,7] = 150.237977760955
$data[10,8] = 38.3760117173083
$data[10,9] = 75.4200443976652
$data[11,0] = 2.53675123794784
$data[11,1] = 52.4431465437837
$data[11,2] = 27.5361165532545
$data[11,3] = 109.831803063784
$data[11,4] = 104.237139366678
$data[11,5] = 96.5307602177052
$data[11,6] = 47.6713096944947
$data[11,7] = 136.000256676227
$data[11,8] = 139.157079597543
$data[11,9] = 28.7473383493476
$data[12,0] = 140.807224223766
$data[12,1] = 52.7106113977314
$data[12,2] = 184.832931302876
$data[12,3] = 67.0119863315541
$data[12,4] = 30.790216862592
$data[12,5] = 176.609724842296
$data[12,6] = 180.237860968447
$data[12,7] = 71.1506246920445
$data[12,8] = 118.938880003495
$data[12,9] = 34.4971641127473
$data[13,0] = 117.622117101039
$data[13,1] = 110.893332357934
$data[13,2] = 79.9020451865634
$data[13,3] = 44.3261203562497
$data[13,4] = 161.165868472851
$data[13,5] = 10.9307847036658
$data[13,6] = 106.114980069043
$data[13,7] = 149.464670172643
$data[13,8] = 3.97442570141257
$data[13,9] = 76.8146358788082
$data[14,0] = 44.2165598479177
$data[14,1] = 162.932115030909
$data[14,2] = 116.769280525283
$data[14,3] = 13.1734127240132
$data[14,4] = 16.2347613909444
$data[14,5] = 157.918465769812
$data[14,6] = 24.8774865758035
$data[14,7] = 194.099975188309
$data[14,8] = 14.9311195197194
$data[14,9] = 179.578932365206
$data[15,0] = 19.8511197324149
$data[15,1] = 49.9856440583177
$data[15,2] = 141.767362943742
$data[15,3] = 151.375606447168
$data[15,4] = 6.87868586130379
$data[15,5] = 144.807912663002
$data[15,6] = 54.6810515479562
$data[15,7] = 187.289033917379
$data[15,8] = 34.1392169865497
$data[15,9] = 192.559849467389
$data[16,0] = 84.9233976029434
$data[16,1] = 14.8735281149268
$data[16,2] = 31.2990977574601
$data[16,3] = 3.87884760456106
$data[16,4] = 157.797927296626
$data[16,5] = 91.6434188800135
$data[16,6] = 172.54110135722
$data[16,7] = 183.209996197005
$data[16,8] = 148.665934590933
$data[16,9] = 93.3063546630118
$data[17,0] = 190.415780148663
$data[17,1] = 98.2066395218515
$data[17,2] = 132.025830974814
$data[17,3] = 62.3424437187344
$data[17,4] = 184.53077850143
$data[17,5] = 177.875109192857
$data[17,6] = 135.941330872449
$data[17,7] = 171.659518578863
$data[17,8] = 50.7772249406098
$data[17,9] = 72.8717510927803
$data[18,0] = 151.732238266492
$data[18,1] = 186.137885780138
$data[18,2] = 56.2195051723251
$data[18,3] = 139.359947638288
$data[18,4] = 14.6460443803324
$data[18,5] = 67.6364730427212
$data[18,6] = 169.125969414192
$data[18,7] = 128.526438739396
$data[18,8] = 37.4474344949459
$data[18,9] = 16.3579558098493
$data[19,0] = 156.24973315571
$data[19,1] = 118.825946151664
$data[19,2] = 115.325453186094
$data[19,3] = 11.4145762340234
$data[19,4] = 191.891238275865
$data[19,5] = 29.3430317329909
$data[19,6] = 131.633017273449
$data[19,7] = 112.890432920722
$data[19,8] = 55.3754854273868
$data[19,9] = 124.591342510931

$ws.Range("A1:J20").Value = $data
